# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For each data row, if the G cell text matches one of the known
# old orderings, rewrite it with the corresponding new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "backup@backdoor.com, System"         = "System, backup@backdoor.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
